$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove cells that are no longer used (content moved / dropped entirely)
# ---------------------------------------------------------------------------
$ws.Range("H3").Clear()
$ws.Range("O5").Clear()
$ws.Range("O6").Clear()
$ws.Range("O7").Clear()
$ws.Range("O8").Clear()
$ws.Range("O9").Clear()
$ws.Range("H10").Clear()
$ws.Range("H11").Clear()

# ---------------------------------------------------------------------------
# Column B ("идея" / need blocks) - shift wording of B7, add new B8
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "нужен прогноз по месяцу"
$ws.Range("B8").Value = " и до конца депозита"

# ---------------------------------------------------------------------------
# New statistics column header, moved to column P
# ---------------------------------------------------------------------------
$ws.Range("P3").Value = "статистика"

# Row 4
$ws.Range("H4").Value = "отчеты"
$ws.Range("H4").Font.Bold = $false
$ws.Range("P4").Value = "DepositExtractor"
$ws.Range("P4").Font.Bold = $true

# Row 5
$ws.Range("H5").Value = "DepositReporter"
$ws.Range("H5").Font.Bold = $true
$ws.Range("P5").Value = "Находит все операции по данному счету"
$ws.Range("P5").Font.Bold = $false

# Row 6
$ws.Range("H6").Value = "составляет List<String> для отчета"
$ws.Range("H6").Font.Bold = $false
$ws.Range("P6").Value = "и составляет таблицу ежедневных остатков"

# Row 7
$ws.Range("H7").Value = "DepositExcelReporter"
$ws.Range("H7").Font.Bold = $true
$ws.Range("P7").Value = "и общие суммы взносов, процентов, расходов"

# Row 8
$ws.Range("H8").Value = "составляет файл экселя"
$ws.Range("K8").Value = "агрегирование"
$ws.Range("M8").Value = "расчет"

# Row 9
$ws.Range("K9").Value = "DepositCalculationAggregator"
$ws.Range("K9").Font.Bold = $true
$ws.Range("M9").Value = "DepositCalculator"
$ws.Range("M9").Font.Bold = $true

# Row 10
$ws.Range("K10").Value = "определяет какой период "
$ws.Range("M10").Value = "расчитывает проценты по вкладу"

# Row 11
$ws.Range("K11").Value = "уже оплачен, какой нет"
$ws.Range("M11").Value = "за каждый день"

# Row 12
$ws.Range("K12").Value = "суммированием определяет "

# Row 13
$ws.Range("K13").Value = "проценты за опред период"

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 25.333333333333332
$ws.Columns.Item(10).ColumnWidth = 17.666666666666668
$ws.Columns.Item(11).ColumnWidth = 27.166666666666668
$ws.Columns.Item(12).ColumnWidth = 9.333333333333334
$ws.Columns.Item(13).ColumnWidth = 21.666666666666668
$ws.Columns.Item(16).ColumnWidth = 17.0

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("H4").Select()
